# Re-colour the solid fills used by the worksheet's coloured-cell palette.
# Each old fill colour (as produced by the original generator) is swapped
# for a new randomly re-rolled colour, while every other attribute of the
# fill / style stays untouched.
#
# Interior.Color in the Excel COM object model is a BGR-ordered long
# (0x00BBGGRR), so the map below is keyed/valued using that encoding of
# the RGB hex values from the workbook's fill table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colorMap = @{
    7664091  = 14726624   # FFDBF174 -> FFE0B5E0
    12089080 = 11985798   # FFF876B8 -> FF86E3B6
    8114055  = 16316318   # FF87CF7B -> FF9EF7F8
    16283228 = 6063450    # FF5C76F8 -> FF5A855C
    15391682 = 11375356   # FFC2DBEA -> FFFC92AD
    8624812  = 9819901    # FFAC9A83 -> FFFDD695
    8609770  = 12935804   # FFEA5F83 -> FF7C62C5
    10350941 = 10745943   # FF5DF19D -> FF57F8A3
    6993874  = 10325682   # FFD2B76A -> FFB28E9D
    14531422 = 14119896   # FF5EBBDD -> FFD873D7
    13891735 = 16750194   # FF97F8D3 -> FF7296FF
    16772462 = 8835685    # FF6EEDFF -> FF65D286
    13346238 = 11464672   # FFBEA5CB -> FFE0EFAE
    6344954  = 15433640   # FFFAD060 -> FFA87FEB
    14272505 = 5755025    # FFF9C7D9 -> FF91D057
    8018789  = 8936280    # FF655B7A -> FF585B88
    15963106 = 6267567    # FFE293F3 -> FFAFA25F
    12320497 = 8877275    # FFF1FEBB -> FFDB7487
    5997498  = 6510736    # FFBA835B -> FF905863
    16474876 = 8388270    # FFFC62FB -> FFAEFE7F
    5699248  = 15881722   # FFB0F656 -> FFFA55F2
    10525910 = 13424574   # FFD69CA0 -> FFBED7CC
    16611501 = 15595006   # FFAD78FD -> FFFEF5ED
    13198018 = 10524273   # FFC262C9 -> FF7196A0
    5727367  = 15781976   # FF876457 -> FF58D0F0
    9888941  = 6051818    # FFADE496 -> FFEA575C
    11436469 = 13534296   # FFB581AE -> FF5884CE
    15638412 = 7581437    # FF8C9FEE -> FFFDAE73
    6464620  = 16341088   # FF6CA462 -> FF6058F9
    9661838  = 14329992   # FF8E6D93 -> FF88A8DA
    12421720 = 9945540    # FF588ABD -> FFC4C197
    6419820  = 5954541    # FF6CF561 -> FFEDDB5A
    11118971 = 12539573   # FF7BA9A9 -> FFB556BF
    14833280 = 16562869   # FF8056E2 -> FFB5BAFC
    12566688 = 6154844    # FFA0C0BF -> FF5CEA5D
}

$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $used.Cells.Item($r, $c)
        $current = $cell.Interior.Color
        if ($colorMap.ContainsKey($current)) {
            $cell.Interior.Color = $colorMap[$current]
        }
    }
}
